$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new row 179 with the latest day's data (2020-08-30) ---
$ws.Range("A179").Value = 44073
$ws.Range("B179").Value = 607938
$ws.Range("C179").Value = 2698
$ws.Range("D179").Value = 79364
$ws.Range("E179").Value = 64007
$ws.Range("F179").Value = 209250
$ws.Range("G179").Value = 25147
$ws.Range("H179").Value = 3550
$ws.Range("I179").Value = 2888
$ws.Range("J179").Value = 5860
$ws.Range("K179").Value = 5247
$ws.Range("L179").Value = 10514
$ws.Range("M179").Value = 3719
$ws.Range("N179").Value = 19568
$ws.Range("O179").Value = 22569
$ws.Range("P179").Value = 5024
$ws.Range("Q179").Value = 4760
$ws.Range("R179").Value = 12289
$ws.Range("S179").Value = 7802
$ws.Range("T179").Value = 13990
$ws.Range("U179").Value = 11608
$ws.Range("V179").Value = 2841
$ws.Range("W179").Value = 1213
$ws.Range("X179").Value = 5872
$ws.Range("Y179").Value = 17724
$ws.Range("Z179").Value = 11624
$ws.Range("AA179").Value = 6952
$ws.Range("AB179").Value = 46549
$ws.Range("AC179").Value = 1029
$ws.Range("AD179").Value = 198
$ws.Range("AE179").Value = 273
$ws.Range("AF179").Value = 446
$ws.Range("AG179").Value = 101
$ws.Range("AH179").Value = 55
$ws.Range("AI179").Value = 265
$ws.Range("AJ179").Value = 1956
$ws.Range("AK179").Value = 3006
$ws.Range("AL179").Value = 36069
$ws.Range("AM179").Value = 6655
$ws.Range("AN179").Value = 2402
$ws.Range("AO179").Value = 36555
$ws.Range("AP179").Value = 914
$ws.Range("AQ179").Value = 20359
$ws.Range("AR179").Value = 1446
$ws.Range("AS179").Value = 8124
$ws.Range("AT179").Value = 1475
$ws.Range("AU179").Value = 1557
$ws.Range("AV179").Value = 4253
$ws.Range("AW179").Value = 1593
$ws.Range("AX179").Value = 934
$ws.Range("AY179").Value = 2466
$ws.Range("AZ179").Value = 2593
$ws.Range("BA179").Value = 46107
$ws.Range("BB179").Value = 12140
$ws.Range("BC179").Value = 2646
$ws.Range("BD179").Value = 7477
$ws.Range("BE179").Value = 3837
$ws.Range("BF179").Value = 278
$ws.Range("BG179").Value = 1395
$ws.Range("BH179").Value = 2584
$ws.Range("BI179").Value = 729
$ws.Range("BJ179").Value = 2013
$ws.Range("BK179").Value = 8308
$ws.Range("BL179").Value = 8322
$ws.Range("BM179").Value = 8365
$ws.Range("BN179").Value = 13763
$ws.Range("BO179").Value = 1877
$ws.Range("BP179").Value = 821
$ws.Range("BQ179").Value = 7483
$ws.Range("BR179").Value = 6517
$ws.Range("BS179").Value = 7608
$ws.Range("BT179").Value = 1500
$ws.Range("BU179").Value = 1501
$ws.Range("BV179").Value = 3003
$ws.Range("BW179").Value = 3157
$ws.Range("BX179").Value = 858
$ws.Range("BY179").Value = 4235
$ws.Range("BZ179").Value = 2413
$ws.Range("CA179").Value = 1293
$ws.Range("CB179").Value = 688
$ws.Range("CC179").Value = 2058
$ws.Range("CD179").Value = 1863
$ws.Range("CE179").Value = 1219
$ws.Range("CF179").Value = 941
$ws.Range("CG179").Value = 4815
$ws.Range("CH179").Value = 1362
$ws.Range("CI179").Value = 1156
$ws.Range("CJ179").Value = 1245
$ws.Range("CK179").Value = 1572
$ws.Range("CL179").Value = 1449
$ws.Range("CM179").Value = 1638
$ws.Range("CN179").Value = 1134
$ws.Range("CO179").Value = 1062
$ws.Range("CP179").Value = 1094
$ws.Range("CQ179").Value = 603
$ws.Range("CR179").Value = 2997
$ws.Range("CS179").Value = 1009
$ws.Range("CT179").Value = 800
$ws.Range("CU179").Value = 735
$ws.Range("CV179").Value = 1274
$ws.Range("CW179").Value = 1165
$ws.Range("CX179").Value = 629
$ws.Range("CY179").Value = 727
$ws.Range("CZ179").Value = 903
$ws.Range("DA179").Value = 1169
$ws.Range("DB179").Value = 975
$ws.Range("DC179").Value = 1110
$ws.Range("DD179").Value = 874
$ws.Range("DE179").Value = 313
$ws.Range("DF179").Value = 332
$ws.Range("DG179").Value = 672
$ws.Range("DH179").Value = 581
$ws.Range("DI179").Value = 403
$ws.Range("DJ179").Value = 530
$ws.Range("DK179").Value = 325
$ws.Range("DL179").Value = 593
$ws.Range("DM179").Value = 700
$ws.Range("DN179").Value = 512
$ws.Range("DO179").Value = 475
$ws.Range("DP179").Value = 363
$ws.Range("DQ179").Value = 512
$ws.Range("DR179").Value = 116981
$ws.Range("DS179").Value = 256604
$ws.Range("DT179").Value = 10246
$ws.Range("DU179").Value = 110546
$ws.Range("DV179").Value = 69505
$ws.Range("DW179").Value = 29381
$ws.Range("DX179").Value = 8802

# Match style of the row above (A column date style, BS:DQ block style)
$ws.Range("A179").Style = $ws.Range("A178").Style
$ws.Range("BS179:DQ179").Style = $ws.Range("BS178:DQ178").Style

# --- Fix up individual data-quality cells: some become "NaN", some get resolved numbers ---
$ws.Range("AB8").Value = "NaN"
$ws.Range("BU13").Value = "NaN"
$ws.Range("BU14").Value = "NaN"
$ws.Range("L18").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("I21").Value = "NaN"
$ws.Range("CM30").Value = "NaN"
$ws.Range("AK33").Value = 1
$ws.Range("AK34").Value = "NaN"
$ws.Range("AK35").Value = "NaN"
$ws.Range("AW62").Value = 2
$ws.Range("AW81").Value = "NaN"
$ws.Range("AP112").Value = "NaN"
$ws.Range("L114").Value = "NaN"
$ws.Range("J115").Value = "NaN"

# --- View state: update frozen-pane anchor and active selection to the new last row ---
$activeWindow = $excel.ActiveWindow
$activeWindow.SplitColumn = 1
$activeWindow.SplitRow = 1
$ws.Range("CZ148").Select()
$activeWindow.FreezePanes = $true
$ws.Range("DR179:DX179").Select()
